# "perbaikan smt dan masa berlaku sertifikat" - add a new "SMT" header
# column to the mahasiswa import-format sheet, between "GEL" and "TAHUN".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column N (the 14th column). This shifts the
# existing "TAHUN" header (previously in N1) and "ID MATERI UJIAN..."
# header (previously in O1) one column to the right, to N1->O1 and
# O1->P1 respectively, opening up N1 for the new column.
$ws.Columns.Item(14).Insert()

# Excel's column-insert recomputes the trailing "default width" column
# block by shifting its range one column to the right; normalize it back
# onto the sheet's valid 1..16384 span by dropping the spurious extra
# column this creates at the very end.
$ws.Columns.Item(16384).Delete()

# Populate the freshly inserted header cell with the new column label.
$ws.Cells.Item(1, 14).Value = "SMT"

# Match the author's view state: scrolled right so column L leads, with
# the new header cell selected.
$ws.Application.ActiveWindow.ScrollColumn = 12
$ws.Range("N1").Select()
